# Insert a new weekly price record as row 58, shifting the existing
# rows 58-97 down to 59-98 (dimension grows from A1:R97 to A1:R98).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 58..97 down by one row, opening up a blank row 58.
$ws.Rows(58).Insert()

# Populate the newly inserted row 58 with the new record's data.
$ws.Range("A58").Value = 1
$ws.Range("B58").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C58").Value = "Arica y Parinacota"
$ws.Range("D58").Value = 45126
$ws.Range("E58").Value = 15
$ws.Range("F58").Value = 100112009
$ws.Range("G58").Value = "Acelga"
$ws.Range("H58").Value = "Sin especificar"
$ws.Range("I58").Value = "Primera"
$ws.Range("J58").Value = 370
$ws.Range("K58").Value = 1800
$ws.Range("L58").Value = 2000
$ws.Range("M58").Value = 1919
$ws.Range("N58").Value = "$/atado 2,5 a 3 kilos"
$ws.Range("O58").Value = "Región de Arica y Parinacota"
$ws.Range("P58").Value = 640
$ws.Range("Q58").Value = 3
$ws.Range("R58").Value = "Hortaliza"
